# Backup QR Scanner data - 4/5/2025, 9:50:36 PM
#
# Adds a new worksheet ("Nghhmsjshd") at the end of the workbook, holding
# the QR-scanner log rows for that location, following the same layout as
# the other "Jhhhsh" / "Kejdhj" / "Hgfg" / "Nghh" backup sheets:
#   Number | Student ID | Location | Log Date | Log Time

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the current last sheet so it lands at the end.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Nghhmsjshd"

# Header row.
$ws.Range("A1").Value = "Number"
$ws.Range("B1").Value = "Student ID"
$ws.Range("C1").Value = "Location"
$ws.Range("D1").Value = "Log Date"
$ws.Range("E1").Value = "Log Time"

# Data rows. Student ID / Log Date / Log Time are stored as text (leading
# apostrophe keeps Excel from reinterpreting the digit/date-like strings as
# numbers or dates), matching how the source data is recorded elsewhere in
# this workbook; "Number" stays a real number.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "'231249"
$ws.Range("C2").Value = "Nghhmsjshd"
$ws.Range("D2").Value = "'2025-04-05"
$ws.Range("E2").Value = "'21:50:18"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "'231249"
$ws.Range("C3").Value = "Nghhmsjshd"
$ws.Range("D3").Value = "'2025-04-05"
$ws.Range("E3").Value = "'21:50:22"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "'231249"
$ws.Range("C4").Value = "Nghhmsjshd"
$ws.Range("D4").Value = "'2025-04-05"
$ws.Range("E4").Value = "'21:50:25"
